$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.097.09"
$ws.Range("E2").Value = "  +3.67%  "
$ws.Range("D3").Value = "2.424.16"
$ws.Range("E3").Value = "  +3.12%  "
$ws.Range("E4").Value = "  +0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "555.01"
$c.ClearFormats()
$ws.Range("E5").Value = "  +2.55%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "138.71"
$c.ClearFormats()
$ws.Range("E6").Value = "  +2.97%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("E9").Value = "  +3.49%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "5.76"
$c.ClearFormats()
$ws.Range("E10").Value = "  +4.05%  "
$ws.Range("E11").Value = "  +0.44%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "24.92"
$c.ClearFormats()
$ws.Range("E13").Value = "  +4.58%  "
$ws.Range("D14").Value = "2.856.83"
$ws.Range("E14").Value = "  +3.14%  "
$ws.Range("D15").Value = "60.030.73"
$ws.Range("E15").Value = "  +3.68%  "
$ws.Range("E16").Value = "  +3.18%  "
$ws.Range("D17").Value = "2.421.55"
$ws.Range("E17").Value = "  +3.27%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "11.38"
$c.ClearFormats()
$ws.Range("E18").Value = "  +6.15%  "
$ws.Range("E19").Value = "  +2.11%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "332.44"
$c.ClearFormats()
$ws.Range("E20").Value = "  +0.71%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.80"
$c.ClearFormats()
$ws.Range("E21").Value = "  +0.94%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E22").Value = "  -0.03%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "65.10"
$c.ClearFormats()
$ws.Range("E23").Value = "  +3.67%  "
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("E25").Value = "  +3.06%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.34"
$c.ClearFormats()
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "0.0₃0789"
$ws.Range("E28").Value = "  +7.18%  "
$ws.Range("E29").Value = "  +1.46%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.29"
$c.ClearFormats()
$ws.Range("E30").Value = "  +2.52%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "168.96"
$c.ClearFormats()
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("E32").Value = "  +2.85%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "18.71"
$c.ClearFormats()
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  +5.09%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.22"
$c.ClearFormats()
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E37").Value = "  +0.07%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.61"
$c.ClearFormats()
$ws.Range("E38").Value = "  +0.26%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.423"
$c.ClearFormats()
$ws.Range("E39").Value = "  +11.91%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "322.01"
$c.ClearFormats()
$ws.Range("E40").Value = "  +11.50%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "39.47"
$c.ClearFormats()
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("E42").Value = "  +1.20%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "140.28"
$c.ClearFormats()
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("E44").Value = "  +1.30%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0521"
$c.ClearFormats()
$ws.Range("E45").Value = "  +2.08%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "19.54"
$c.ClearFormats()
$ws.Range("E46").Value = "  +2.10%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.415"
$c.ClearFormats()
$ws.Range("E47").Value = "  +8.75%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.575"
$c.ClearFormats()
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("E49").Value = "  +2.01%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "17.81"
$c.ClearFormats()
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("E51").Value = "  -0.18%  "
